$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.992.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.68%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.862.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.86%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.47%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4697"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3900"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.83"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.82%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07966"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9795"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.69%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.51"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.19%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.932"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.39%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.843.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.197"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.50%  "

$ws.Range("E17").Value = "  -0.42%  "

$ws.Range("E18").Value = "  +0.48%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06623"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.83%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.985.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.61%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.391"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.74%  "

$ws.Range("E25").Value = "  -0.88%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.42%  "

$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.067.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.43%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.28%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.098"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.66%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.463"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.46%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9578"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.93%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09480"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.90%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.578"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.40%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.305"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.346"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06082"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.67%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02246"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.81%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.286"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.164"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.001"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.30%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5912"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.34%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1865"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.21%  "

$ws.Range("E44").Value = "  -1.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.281"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5539"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.82%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.948"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.97%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06877"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.71%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.97%  "

$ws.Range("E51").Value = "  -32.78%  "
